$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$checkMark = [string][char]0x2705

$oldLine1 = $checkMark + " 1000 Bs = 3.25 = 12460.99 pesos"
$newLine1 = $checkMark + " 1000 Bs = 3.25 = 12426.54 pesos"

$oldLine2 = $checkMark + " 12460.99 pesos = 3.24 = 971.27 Bs"
$newLine2 = $checkMark + " 12426.54 pesos = 3.23 = 969.47 Bs"

$text = $wsHoja1.Range("A1").Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update numeric rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 307.97
$wsTasas.Range("O10").Value = 3827
$wsTasas.Range("N12").Value = 3845.48
